$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (same style as existing header cells, e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# New data values for columns I and J (rows 2-10)
$values = @{
    2  = @(1, 5)
    3  = @(8, 8)
    4  = @(1, 2)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 3)
    8  = @(1, 4)
    9  = @(1, 3)
    10 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
